# Update hashcode values in the "hashcode.csv" sheet (column B) to reflect
# the automatic hashcode regeneration described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hashcode.csv")

$updates = @{
    "B9"   = "a6037f74bda70b9c0578221cc3123f03"
    "B11"  = "cc3b2caf8f96facc53c369be93550474"
    "B34"  = "a7dd9997e1c40f89a8f7dbfc9c891cad"
    "B162" = "46a143b70833780e4cf4fcd8422cb6c2"
    "B180" = "0a6515b9b92c973e3013b6f12968668c"
    "B191" = "5cb153fcce5c7fc167b7711dd0b9e59d"
    "B213" = "d1bc28f873fd85d7c32d58c113ae62d7"
    "B227" = "eacb84044618ca5941a26f8bb17104c0"
    "B232" = "f626a3db4a98dbcf6ac08a1230606469"
    "B419" = "2ee5add6736bc97726d8045230c25adb"
    "B461" = "d630e0a02237c5cb7cc8fdacad527d79"
    "B478" = "911324c32a26fbe7007b2e2ebaef8187"
    "B506" = "25f2e490b376c79ab16e17eb5089138e"
    "B514" = "453191517d23d7051c0d303be15ae0b4"
    "B524" = "408e02ddc0b6f4215af6415b3ddf22a7"
    "B666" = "429c930454134a12c1592d5829630dd0"
    "B680" = "7a7b55ebde53e22a1b4e021f6bc4ff84"
    "B685" = "160e0e1b3c18eb934ad87655e4edcf22"
    "B703" = "fe9579abfdc3819b5efb7c1d99932e09"
    "B729" = "322310ae54f161a29946ec86c092b066"
    "B854" = "07bb12e1a04cece46fbb9f865931cd62"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
